{"js": "// Update the date heading and the practice-table answers to the next day's\n// values. Cell text is replaced via each paragraph's Range so the existing\n// run/paragraph formatting (fonts, size, alignment) is preserved instead of\n// being reset by a body-level insertText.\n\n// 1) Update the date paragraph (first paragraph in the body).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateRange = paragraphs.items[0].getRange();\ndateRange.insertText(\"2024-06-23 Sunday\", Word.InsertLocation.replace);\n\n// 2) Update the answers table. The table has 5 columns and a data row\n// every 4th row (rows 0, 4, 8, 12, 16), the rows in-between being blank\n// placeholders for handwritten work. Only the cell text content changes;\n// one cell (row 4, col 0) keeps its original value.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst newValues = [\n  [0, 0, \"41\u00f74=10, 1\"],\n  [0, 1, \"86\u00f75=17, 1\"],\n  [0, 2, \"76\u00f76=12, 4\"],\n  [0, 3, \"50\u00f76=8, 2\"],\n  [0, 4, \"33\u00f74=8, 1\"],\n  [4, 1, \"16\u00f76=2, 4\"],\n  [4, 2, \"69\u00f73=23, 0\"],\n  [4, 3, \"27\u00f77=3, 6\"],\n  [4, 4, \"30\u00f79=3, 3\"],\n  [8, 0, \"29\u00f76=4, 5\"],\n  [8, 1, \"39\u00f79=4, 3\"],\n  [8, 2, \"16\u00f78=2, 0\"],\n  [8, 3, \"43\u00f74=10, 3\"],\n  [8, 4, \"21\u00f76=3, 3\"],\n  [12, 0, \"20\u00f73=6, 2\"],\n  [12, 1, \"69\u00f78=8, 5\"],\n  [12, 2, \"61\u00f72=30, 1\"],\n  [12, 3, \"50\u00f73=16, 2\"],\n  [12, 4, \"17\u00f72=8, 1\"],\n  [16, 0, \"42\u00f78=5, 2\"],\n  [16, 1, \"56\u00f73=18, 2\"],\n  [16, 2, \"32\u00f74=8, 0\"],\n  [16, 3, \"65\u00f78=8, 1\"],\n  [16, 4, \"40\u00f78=5, 0\"],\n];\n\nfor (const [row, col, text] of newValues) {\n  const cell = table.getCell(row, col);\n  const cellRange = cell.body.paragraphs.getFirst().getRange();\n  cellRange.insertText(text, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and the practice-table answers to the next day's\n# values. Assigning to Range.Text replaces only the run's text content, so\n# the existing paragraph/run formatting (fonts, size, alignment) is kept.\n\n$d = $word.ActiveDocument\n\n# 1) Update the date paragraph (first paragraph in the body).\n$d.Paragraphs.Item(1).Range.Text = \"2024-06-23 Sunday\"\n\n# 2) Update the answers table. The table has 5 columns and a data row\n# every 4th row (rows 1, 5, 9, 13, 17 in 1-based COM indexing), the rows\n# in-between being blank placeholders for handwritten work. Only the cell\n# text content changes; one cell (row 5, col 1) keeps its original value.\n$t = $d.Tables.Item(1)\n\n$t.Cell(1, 1).Range.Text = \"41\u00f74=10, 1\"\n$t.Cell(1, 2).Range.Text = \"86\u00f75=17, 1\"\n$t.Cell(1, 3).Range.Text = \"76\u00f76=12, 4\"\n$t.Cell(1, 4).Range.Text = \"50\u00f76=8, 2\"\n$t.Cell(1, 5).Range.Text = \"33\u00f74=8, 1\"\n\n$t.Cell(5, 2).Range.Text = \"16\u00f76=2, 4\"\n$t.Cell(5, 3).Range.Text = \"69\u00f73=23, 0\"\n$t.Cell(5, 4).Range.Text = \"27\u00f77=3, 6\"\n$t.Cell(5, 5).Range.Text = \"30\u00f79=3, 3\"\n\n$t.Cell(9, 1).Range.Text = \"29\u00f76=4, 5\"\n$t.Cell(9, 2).Range.Text = \"39\u00f79=4, 3\"\n$t.Cell(9, 3).Range.Text = \"16\u00f78=2, 0\"\n$t.Cell(9, 4).Range.Text = \"43\u00f74=10, 3\"\n$t.Cell(9, 5).Range.Text = \"21\u00f76=3, 3\"\n\n$t.Cell(13, 1).Range.Text = \"20\u00f73=6, 2\"\n$t.Cell(13, 2).Range.Text = \"69\u00f78=8, 5\"\n$t.Cell(13, 3).Range.Text = \"61\u00f72=30, 1\"\n$t.Cell(13, 4).Range.Text = \"50\u00f73=16, 2\"\n$t.Cell(13, 5).Range.Text = \"17\u00f72=8, 1\"\n\n$t.Cell(17, 1).Range.Text = \"42\u00f78=5, 2\"\n$t.Cell(17, 2).Range.Text = \"56\u00f73=18, 2\"\n$t.Cell(17, 3).Range.Text = \"32\u00f74=8, 0\"\n$t.Cell(17, 4).Range.Text = \"65\u00f78=8, 1\"\n$t.Cell(17, 5).Range.Text = \"40\u00f78=5, 0\"\n"}
